# Bugfix in Edge Weighting Script
# Re-point the sheet's AutoFilter from column A ("Weight" == 0.1) to
# column B ("Outgoing" == the selected course), which changes which rows
# are hidden, and leaves the selection on the cell that was used to pick
# the new filter value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing AutoFilter criteria (column A filtered on 0.1 / "0,001")
$ws.AutoFilter.Range.AutoFilter()

# Re-apply the AutoFilter over the same range, now filtering column B
# (field 2) to the single course value, using xlFilterValues (7) so the
# criteria is written out as a plain <filters><filter val="..."/></filters>
# list rather than a <customFilter>.
$ws.Range("A1:E91").AutoFilter(2, @("CS 302 Praktische Informatik I Practical Computer Science I"), 7)

# Selection moved to the cell whose value was used to build the filter.
$ws.Range("B22").Select()
